$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark a couple of items as resolved / closed
$ws.Range("B24").Value = "Resolvido"
$ws.Range("B26").Value = "-"

# Add a note in column D explaining that the table is in the article
$ws.Range("D26").Value = "A tabela está no artigo"

# Move the "A Tabela 1 não é citada no texto!" row (row 28) above the
# "Na hipótese H0..." row (row 27), and assign it to Laura.
$ws.Range("A27").Value = "A Tabela 1 não é citada no texto!"
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "Laura"

$ws.Range("A28").Value = "Na hipótese H0 o que são resultados insuficientes?"
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = ""

# Mark two more items as resolved
$ws.Range("B38").Value = "Resolvido"
$ws.Range("B41").Value = "Resolvido"

# Update the visible window / selection to reflect where the user ended up
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("A29").Select()
